$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 (value 685601396368), shifting all subsequent rows up by one.
$ws.Rows.Item(2).Delete()

# Update selection to match the target state (single cell K5).
$ws.Range("K5").Select()
